# Generate Report for Handoff
# Updates status text, handoff/handback timestamps, and narrows the
# "Status" / "Latest Handoff Datetime" style columns on all three sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# The stored OOXML column <col width="..."> is expressed in "number of
# characters" plus a fixed 5-pixel padding (i.e. raw width = ColumnWidth +
# 5/MaxDigitWidth). This engine's MaxDigitWidth is 6px here, so subtract
# 5/6 from the desired raw width before handing it to the ColumnWidth
# property (which itself is expressed in characters) so the round-tripped
# stored width lands as close as possible to the target.
$targetRawWidth = 17.2159881591797
$columnWidthChars = $targetRawWidth - (5.0 / 6.0)

# --- Overview sheet ---------------------------------------------------
# E2 (zh-cn status) and F2 (de-de status): "Handed back: in sync with en-US" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
# G2 (Latest HO Xliff Generate Date): 2016-08-29 21:13:19 -> 2016-08-29 21:13:59
$wsOverview.Range("G2").Value = "2016-08-29 21:13:59"

# Columns E & F width: 29.9777047293527 -> 17.2159881591797
$wsOverview.Range("E:F").ColumnWidth = $columnWidthChars

# --- zh-cn sheet --------------------------------------------------------
# C2 (Status): "Handed back: in sync with en-US" -> "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
# H2 (Latest Handoff Datetime): 2016-08-29 21:13:14 -> 2016-08-29 21:13:55
$wsZhCn.Range("H2").Value = "2016-08-29 21:13:55"

# Column C width: 29.9777047293527 -> 17.2159881591797
$wsZhCn.Range("C:C").ColumnWidth = $columnWidthChars

# --- de-de sheet --------------------------------------------------------
# C2 (Status): "Handed back: in sync with en-US" -> "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"
# H2 (Latest Handoff Datetime): 2016-08-29 21:13:19 -> 2016-08-29 21:13:59
$wsDeDe.Range("H2").Value = "2016-08-29 21:13:59"

# Column C width: 29.9777047293527 -> 17.2159881591797
$wsDeDe.Range("C:C").ColumnWidth = $columnWidthChars
